$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "26.677.56"
Set-TextValue "E2" "  -1.33%  "

Set-TextValue "D3" "1.594.31"
Set-TextValue "E3" "  -1.48%  "

Set-TextValue "E4" "  -0.14%  "

Set-TextValue "D5" "211.52"
Set-TextValue "E5" "  -1.23%  "

Set-TextValue "D6" "0.512"
Set-TextValue "E6" "  -0.79%  "

Set-TextValue "E7" "  -0.12%  "

Set-TextValue "D8" "0.247"
Set-TextValue "E8" "  -2.04%  "

Set-TextValue "D9" "0.0617"
Set-TextValue "E9" "  -1.43%  "

Set-TextValue "D10" "19.60"
Set-TextValue "E10" "  -2.06%  "

Set-TextValue "D11" "0.0834"
Set-TextValue "E11" "  -1.39%  "

Set-TextValue "D12" "1.815.47"
Set-TextValue "E12" "  -1.64%  "

Set-TextValue "D13" "1.589.28"
Set-TextValue "E13" "  -1.93%  "

Set-TextValue "D14" "4.04"
Set-TextValue "E14" "  -2.56%  "

Set-TextValue "D15" "0.523"
Set-TextValue "E15" "  -3.01%  "

Set-TextValue "D16" "65.21"
Set-TextValue "E16" "  +0.92%  "

Set-TextValue "D17" "26.658.87"
Set-TextValue "E17" "  -1.39%  "

Set-TextValue "D18" "0.0₃0729"
Set-TextValue "E18" "  -2.59%  "

Set-TextValue "B19" "BitcoinCash"
Set-TextValue "C19" "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue "D19" "208.53"
Set-TextValue "E19" "  -2.50%  "

Set-TextValue "B20" "Dai"
Set-TextValue "C20" "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue "D20" "1.00"
Set-TextValue "E20" "  -0.10%  "

Set-TextValue "D21" "6.72"
Set-TextValue "E21" "  -1.75%  "

Set-TextValue "D22" "4.27"
Set-TextValue "E22" "  -2.18%  "

Set-TextValue "D23" "2.33"
Set-TextValue "E23" "  -3.38%  "

Set-TextValue "D24" "8.91"
Set-TextValue "E24" "  -1.40%  "

Set-TextValue "D25" "145.97"
Set-TextValue "E25" "  -1.39%  "

Set-TextValue "E26" "  -0.02%  "

Set-TextValue "D27" "7.19"
Set-TextValue "E27" "  -3.07%  "

Set-TextValue "E28" "  -1.96%  "

Set-TextValue "D29" "15.33"
Set-TextValue "E29" "  -1.32%  "

Set-TextValue "D30" "0.0504"
Set-TextValue "E30" "  -2.45%  "

Set-TextValue "E31" "  -1.77%  "

Set-TextValue "D32" "3.24"
Set-TextValue "E32" "  -3.81%  "

Set-TextValue "D33" "0.671"
Set-TextValue "E33" "  -15.09%  "

Set-TextValue "D34" "2.91"
Set-TextValue "E34" "  -3.00%  "

Set-TextValue "D35" "1.301.56"
Set-TextValue "E35" "  -2.64%  "

Set-TextValue "D36" "2.44"
Set-TextValue "E36" "  -0.70%  "

Set-TextValue "E37" "  -4.61%  "

Set-TextValue "E38" "  -4.02%  "

Set-TextValue "D39" "0.830"
Set-TextValue "E39" "  -2.99%  "

Set-TextValue "E40" "  -0.05%  "

Set-TextValue "D41" "0.791"
Set-TextValue "E41" "  -1.43%  "

Set-TextValue "D42" "5.35"
Set-TextValue "E42" "  +1.23%  "

Set-TextValue "D43" "2.18"
Set-TextValue "E43" "  -1.98%  "

Set-TextValue "D44" "63.17"
Set-TextValue "E44" "  -3.36%  "

Set-TextValue "D45" "1.728.90"
Set-TextValue "E45" "  -1.36%  "

Set-TextValue "D46" "89.22"
Set-TextValue "E46" "  -0.58%  "

Set-TextValue "E47" "  -2.25%  "

Set-TextValue "D48" "0.824"
Set-TextValue "E48" "  -8.70%  "

Set-TextValue "D49" "0.0503"
Set-TextValue "E49" "  -2.54%  "

Set-TextValue "B50" "Algorand"
Set-TextValue "C50" "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue "D50" "0.0974"
Set-TextValue "E50" "  -2.58%  "

Set-TextValue "B51" "EnergySwap"
Set-TextValue "C51" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D51" "7.50"
Set-TextValue "E51" "  -1.58%  "
